$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet ---
$ws.Name = "Corrected data"

# --- Header text ---
$ws.Range("B1").Value = "people"

# --- New data values (column B, "people") ---
$ws.Range("B2").Value = 153706850
$ws.Range("B3").Value = 159061181
$ws.Range("B4").Value = 166151041
$ws.Range("B5").Value = 173260603
$ws.Range("B6").Value = 181281186
$ws.Range("B7").Value = 185155513
$ws.Range("B8").Value = 187965778
$ws.Range("B9").Value = 192913686
$ws.Range("B10").Value = 188228921

# --- Tab color ---
$ws.Tab.ThemeColor = 6
$ws.Tab.TintAndShade = 0.39997558519241921

# --- Selection ---
$ws.Range("B2").Select()

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 10.83203125
$ws.Columns("B").ColumnWidth = 32.5
$ws.Columns("C").ColumnWidth = 25.5

# --- Remove the shaded header fill (match it to the plain body style first) ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Body font size 12 -> 15 (years + values) ---
$ws.Range("A1:B10").Font.Size = 15

# --- Number format (thousands separator) on the "people" column, and header row ---
$ws.Range("B2:B10").NumberFormat = "#,##0"
$ws.Range("A1:B1").NumberFormat = "#,##0"

# --- Header styling: bold, size 14, vertically centered ---
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Size = 14
$ws.Range("A1:B1").Font.ThemeColor = 1
$ws.Range("A1:B1").VerticalAlignment = -4108

# --- Row heights ---
$ws.Rows(1).RowHeight = 18
$ws.Range("A2:A10").EntireRow.RowHeight = 19

Write-Host "edit applied"
